$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated QUANTITY (D) and AVERAGE PRICE (E) figures for each holding.
# Also apply a numeric format (0.00) to the AVERAGE PRICE column, which
# introduces a new cell style used by E2:E31.

$ws.Range("D2").Value = 6468
$ws.Range("E2").Value = 28.740374149659864

$ws.Range("D3").Value = 2987
$ws.Range("E3").Value = 61.196317375292935

$ws.Range("D4").Value = 5159
$ws.Range("E4").Value = 93.686865671641797

$ws.Range("D5").Value = 660
$ws.Range("E5").Value = 128.60936363636364

$ws.Range("D6").Value = 1005
$ws.Range("E6").Value = 91.778208955223874

$ws.Range("D7").Value = 18201
$ws.Range("E7").Value = 24.297907807263339

$ws.Range("D8").Value = 16826
$ws.Range("E8").Value = 31.489545346487578

$ws.Range("D9").Value = 11228
$ws.Range("E9").Value = 32.797620235126466

$ws.Range("D10").Value = 3475
$ws.Range("E10").Value = 59.316644604316544

$ws.Range("D11").Value = 761
$ws.Range("E11").Value = 110.25467805519054

$ws.Range("D12").Value = 454
$ws.Range("E12").Value = 153.05731277533042

$ws.Range("D13").Value = 3449
$ws.Range("E13").Value = 99.014650623369079

$ws.Range("D14").Value = 4854
$ws.Range("E14").Value = 41.215852904820764

$ws.Range("D15").Value = 1624
$ws.Range("E15").Value = 155.48120689655173

$ws.Range("D16").Value = 428
$ws.Range("E16").Value = 71.173738317757

$ws.Range("D17").Value = 13676
$ws.Range("E17").Value = 10.399371161158232

$ws.Range("D18").Value = 19951
$ws.Range("E18").Value = 23.099100796952534

$ws.Range("D19").Value = 2445
$ws.Range("E19").Value = 31.200507157464216

$ws.Range("D20").Value = 1057
$ws.Range("E20").Value = 88.056944181646159

$ws.Range("D21").Value = 504
$ws.Range("E21").Value = 233.71541666666667

$ws.Range("D22").Value = 16428
$ws.Range("E22").Value = 24.398210372534692

$ws.Range("D23").Value = 1108
$ws.Range("E23").Value = 89.49648014440433

$ws.Range("D24").Value = 5932
$ws.Range("E24").Value = 16.858755900202294

$ws.Range("D25").Value = 1362
$ws.Range("E25").Value = 294.59596916299563

$ws.Range("D26").Value = 12969
$ws.Range("E26").Value = 12.209190377052973

$ws.Range("D27").Value = 4810
$ws.Range("E27").Value = 86.491349272349268

$ws.Range("D28").Value = 6752
$ws.Range("E28").Value = 28.978667061611372

$ws.Range("D29").Value = 482
$ws.Range("E29").Value = 978.51375518672205

$ws.Range("D30").Value = 469
$ws.Range("E30").Value = 172.0038805970149

$ws.Range("D31").Value = 805
$ws.Range("E31").Value = 98.433105590062127

# Apply a 2-decimal numeric display format to the AVERAGE PRICE column,
# which creates the new cellXfs entry (numFmtId 2) seen in the diff.
$ws.Range("E2:E31").NumberFormat = "0.00"

# Move the active selection to match the author's last cursor position.
$null = $ws.Range("F16").Select()
